# edit.ps1 - apply the CV edits described by the unified diff
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $null = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# 1. PRECEPTRON-PTMKB -> PERCEPTRON-PTMKB (typo fix, bold run)
Replace-Text "PRECEPTRON-PTMKB" "PERCEPTRON-PTMKB"

# 2. Reword the PTMKB project description bullet
Replace-Text `
    "Experimenting with a novel deep learning model for residue-level PTM prediction by considering structural conformation and pairwise relationships between residues" `
    "Developing a deep learning model for residue-level PTM prediction using structural conformation and residue-pair relationships"

# 3. Extend the oncology / clinical-notes bullet with a new clause
Replace-Text `
    "Evaluated role of language models for processing oncology-based clinical notes and extracting phenotypes" `
    "Evaluated role of language models for processing oncology-based clinical notes and extracting phenotypes; compared against ontology-based baselines"

# 4. BDI architecture bullet: "LLMs" -> "language model planners"
Replace-Text `
    "Proposed a novel agentic architecture utilizing the Belief-Desire-Intention (BDI) model integrating LLMs with BDI agents" `
    "Proposed a novel agentic architecture utilizing the Belief-Desire-Intention (BDI) model integrating language model planners with BDI agents"

# 5. Digital twin bullet: promote from sub-bullet (ilvl 1, indent 1440) to top-level
#    bullet (ilvl 0, indent 720), plus reword the text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*digital twin system prototype*") {
        $p.Range.ListFormat.ListLevelNumber = 1
        $p.Range.ParagraphFormat.LeftIndent = 36
    }
}
Replace-Text `
    "Implemented a digital twin system prototype achieving hundreds of thousands of agents running simultaneously" `
    "Implemented a digital twin system prototype successfully scaled to >100k simulated agents in stress tests on a 4-node GPU cluster"

# 6. PERCEPTRON-PTMKB publication line: drop the journal parenthetical and
#    turn the status note from italic "Submitted Oct 2025, Under Review"
#    into bold "Submitted Oct 2025; under review"
Replace-Text " (Computational Resources for Molecular Biology), " ", "
$statusRng = $d.Content
if ($statusRng.Find.Execute("Submitted Oct 2025, Under Review")) {
    $statusRng.Text = "Submitted Oct 2025; under review"
    $statusRng.Italic = 0
    $statusRng.Bold = 1
}

# 7. HAKI publication line: italic "In Preparation" -> bold "Manuscript In Preparation"
$prepRng = $d.Content
if ($prepRng.Find.Execute("In Preparation")) {
    $prepRng.Text = "Manuscript In Preparation"
    $prepRng.Italic = 0
    $prepRng.Bold = 1
}

# 8. Novozymes project: merge the two sub-bullets ("Applied feature engineering
#    on limited dataset" and "Achieved a Spearman's correlation of 0.56") into
#    the parent bullet's second run, then delete the now-empty paragraphs.
Replace-Text `
    " Developed an XGBoost-based forest model for finding the optimal melting temperature of a trial enzyme" `
    " Developed an XGBoost-based forest model for finding the optimal melting temperature of a trial enzyme; Spearman’s ρ = 0.56 on limited data"

$featureEngPara = $null
$spearmanPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*feature engineering*") { $featureEngPara = $p }
    if ($p.Range.Text -like "*Spearman*correlation*") { $spearmanPara = $p }
}
if ($spearmanPara -ne $null) { $spearmanPara.Range.Delete() }
if ($featureEngPara -ne $null) { $featureEngPara.Range.Delete() }

# 9. TORCS project bullet rewording
Replace-Text `
    " Designed an artificial neural network to make a car drive intelligently based on the features captured through hundreds of hours of manual driving" `
    " Trained a racing-policy ANN in TORCS from hundreds of hours of human driving logs"

# 10. Skills heading: "Programming Language:" -> "Programming Languages:"
Replace-Text "Programming Language:" "Programming Languages:"

Write-Output "All edits applied."
